$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "Ativação:" date value, 01/01/2016 -> 01/01/2023 ---
# The literal string "01/01/2023" looks like a date to Excel's Value setter,
# which would auto-convert it to a date serial and mint a brand new
# number-formatted style. To keep it as plain text (as in the source file)
# we stage it in a scratch cell that's explicitly formatted as Text first,
# then copy only the (already-text) VALUE into B8/C8, and finally restore
# B8/C8's original cell formatting (wrap text, vertical top, font) by
# copying FORMATS from the untouched B9/C9 pair, which already carries the
# exact formatting B8/C8 should keep.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"
$scratch.Value = "01/01/2023"

$scratch.Copy()
$ws.Range("B8").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("B9").Copy()
$ws.Range("B8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C9").Copy()
$ws.Range("C8").PasteSpecial(-4122)   # xlPasteFormats

$scratch.Clear()

# --- Row 11: "Objectives:" - add English objectives text in B/C ---
$objectivesText = "Provide the student with the basic knowledge of electronic materials aiming their application in devices."
$ws.Range("B11").Value = $objectivesText
$ws.Range("C11").Value = $objectivesText
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 14: "Short syllabus:" - add English short syllabus text in B/C ---
$shortSyllabusText = "Materials for electronics. Electronics and Solid State Physics. Semiconductor materials and devices. Optoelectronic materials and devices. Dielectric and piezoelectric materials and devices."
$ws.Range("B14").Value = $shortSyllabusText
$ws.Range("C14").Value = $shortSyllabusText
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 16: "Syllabus:" - add English full syllabus text in B/C ---
$syllabusText = "Materials for electronic applications: metals, ceramics, glasses and polymers. Single crystals and thin films.Waves and particles in matter. Electrons in atoms and crystals. Energy band structures. Electronic and spectroscopic properties of materials.Conducting, semiconducting and insulating materials. Electronic properties in semiconductors. Electric transport. Semiconductor devices. pn junction Metal-semiconductor and semiconductor-insulator contact. Semiconductor devices: diodes and bipolar and FET transistors.Optoelectronic materials and devices. LED, semiconductor laser, photodetectors and photovoltaic cells.Types and properties of dielectric materials. Ferroelectric and piezoelectric materials. Devices based on dielectric and piezoelectric materials. Applications."
$ws.Range("B16").Value = $syllabusText
$ws.Range("C16").Value = $syllabusText
$ws.Range("B13").Copy()
$ws.Range("B16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C13").Copy()
$ws.Range("C16").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
